# HCAP-1233: update participants_status.xlsx seed-data fixture
#  - the test "site" id used by the seed-data JSON blob moved from 4 -> 1
#  - the participant_id sample value moved from 68 -> 26
#  - leave the cursor on H5, matching where the author left off editing

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# participant_id sample value
$ws.Range("C2").Value = 26

# data column: seed "site" changed from 4 to 1 (rest of the JSON payload unchanged)
$ws.Range("G2").Value = '{"site":1,"hiredDate":"2022/05/02","startDate":"2022/05/03","positionType":"","positionTitle":"","nonHcapOpportunity":false}'

# leave selection where the author left it
[void]$ws.Range("H5").Select()
